$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("test1")
$ws1.Name = "testValidLogin"

$ws1.Cells.Item(1,1).Value = "UN"
$ws1.Cells.Item(1,2).Value = "PW"
$ws1.Cells.Item(1,3).Value = "E-Title"
$ws1.Cells.Item(2,2).Value = "manager"
$ws1.Cells.Item(2,3).Value = "actiTIME - Enter Time-Track"
$ws1.Cells.Item(2,1).Value = "ADMIN"

$ws1.Rows.Item(4).Delete()
$ws1.Rows.Item(3).Delete()

$ws1.Columns.Item(3).EntireColumn.AutoFit() | Out-Null

$ws1.Select()
$ws1.Range("A2").Select()
